$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Range, $Value)
    $Range.NumberFormat = "@"
    $Range.Value = $Value
    $Range.NumberFormat = "General"
}

$ws.Range("D2").Value = '45.297.61'
$ws.Range("E2").Value = '  -3.24%  '
$ws.Range("D3").Value = '2.426.96'
$ws.Range("E3").Value = '  +7.19%  '
Set-TextValue $ws.Range("D4") '1.00'
$ws.Range("E4").Value = '  -0.01%  '
Set-TextValue $ws.Range("D5") '293.53'
$ws.Range("E5").Value = '  -2.44%  '
Set-TextValue $ws.Range("D6") '93.81'
$ws.Range("E6").Value = '  -6.27%  '
$ws.Range("E7").Value = '  -0.52%  '
Set-TextValue $ws.Range("D8") '1.00'
$ws.Range("E8").Value = '  +0.02%  '
$ws.Range("E10").Value = '  -3.93%  '
$ws.Range("E11").Value = '  -0.45%  '
$ws.Range("E12").Value = '  -2.13%  '
$ws.Range("E13").Value = '  +1.83%  '
$ws.Range("D14").Value = '2.798.58'
$ws.Range("E14").Value = '  +7.42%  '
$ws.Range("D15").Value = '2.435.55'
$ws.Range("E15").Value = '  +7.74%  '
Set-TextValue $ws.Range("D16") '14.22'
$ws.Range("E16").Value = '  +4.97%  '
$ws.Range("E17").Value = '  +5.19%  '
$ws.Range("D18").Value = '45.248.32'
$ws.Range("E18").Value = '  -3.28%  '
Set-TextValue $ws.Range("D19") '12.38'
$ws.Range("E19").Value = '  -2.38%  '
$ws.Range("D20").Value = '0.0₃0938'
$ws.Range("E20").Value = '  +1.39%  '
Set-TextValue $ws.Range("D21") '6.19'
$ws.Range("E21").Value = '  +5.67%  '
Set-TextValue $ws.Range("D22") '67.07'
$ws.Range("E22").Value = '  +3.20%  '
Set-TextValue $ws.Range("D23") '239.26'
$ws.Range("E23").Value = '  -3.67%  '
Set-TextValue $ws.Range("D24") '2.76'
$ws.Range("E24").Value = '  -1.41%  '
Set-TextValue $ws.Range("D25") '1.00'
$ws.Range("E25").Value = '  -0.01%  '
Set-TextValue $ws.Range("D26") '1.91'
$ws.Range("E26").Value = '  +2.49%  '
$ws.Range("E27").Value = '  -0.86%  '
Set-TextValue $ws.Range("D28") '37.18'
$ws.Range("E28").Value = '  -12.41%  '
Set-TextValue $ws.Range("D29") '9.55'
$ws.Range("E29").Value = '  -1.64%  '
Set-TextValue $ws.Range("D30") '3.86'
$ws.Range("E30").Value = '  +20.86%  '
$ws.Range("E31").Value = '  +7.54%  '
Set-TextValue $ws.Range("D32") '149.22'
$ws.Range("E32").Value = '  +2.60%  '
Set-TextValue $ws.Range("D33") '2.72'
$ws.Range("E33").Value = '  -2.27%  '
Set-TextValue $ws.Range("D34") '5.40'
$ws.Range("E34").Value = '  +0.11%  '
$ws.Range("E35").Value = '  -1.62%  '
$ws.Range("E36").Value = '  +16.79%  '
$ws.Range("E37").Value = '  -1.49%  '
$ws.Range("E38").Value = '  -0.61%  '
Set-TextValue $ws.Range("D39") '14.29'
$ws.Range("E39").Value = '  -11.59%  '
Set-TextValue $ws.Range("D40") '3.72'
$ws.Range("E40").Value = '  -2.33%  '
Set-TextValue $ws.Range("D41") '0.0294'
$ws.Range("E41").Value = '  -1.25%  '
$ws.Range("D42").Value = '1.993.10'
$ws.Range("E42").Value = '  +11.45%  '
Set-TextValue $ws.Range("D43") '3.16'
$ws.Range("E43").Value = '  -1.26%  '
$ws.Range("E44").Value = '  +0.11%  '
Set-TextValue $ws.Range("D45") '88.39'
$ws.Range("E45").Value = '  -2.86%  '
Set-TextValue $ws.Range("D46") '16.20'
$ws.Range("E46").Value = '  +27.47%  '
Set-TextValue $ws.Range("D47") '1.70'
$ws.Range("E47").Value = '  -13.10%  '
$ws.Range("E48").Value = '  +9.82%  '
Set-TextValue $ws.Range("D49") '101.14'
$ws.Range("E49").Value = '  +7.90%  '
$ws.Range("D50").Value = '2.666.89'
$ws.Range("E50").Value = '  +7.38%  '
$ws.Range("E51").Value = '  -3.68%  '
